$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# New header labels in columns L and M of Sheet1
$ws1.Range("L1").Value = "Center-to-center distance"
$ws1.Range("M1").Value = "Pillar diameter"

# New column widths (closest achievable match to the bestFit widths of the target sheet;
# the headless engine quantizes ColumnWidth to 1/6-character steps)
$ws1.Columns.Item(12).ColumnWidth = 23.5
$ws1.Columns.Item(13).ColumnWidth = 13.5

# New data values for rows 2-11 (Center-to-center distance, Pillar diameter)
$ws1.Range("L2").Value = 30
$ws1.Range("M2").Value = 16

$ws1.Range("L3").Value = 30
$ws1.Range("M3").Value = 16

$ws1.Range("L4").Value = 25
$ws1.Range("M4").Value = 16

$ws1.Range("L5").Value = 30
$ws1.Range("M5").Value = 20

$ws1.Range("L6").Value = 30
$ws1.Range("M6").Value = 20

$ws1.Range("L7").Value = 30
$ws1.Range("M7").Value = 20

$ws1.Range("L8").Value = 35
$ws1.Range("M8").Value = 20

$ws1.Range("L9").Value = 35
$ws1.Range("M9").Value = 20

$ws1.Range("L10").Value = 35
$ws1.Range("M10").Value = 20

$ws1.Range("L11").Value = 35
$ws1.Range("M11").Value = 20

# Sheet2's selection stays at E9 (it just loses the "active tab" state below)
$ws2.Activate()
$ws2.Range("E9").Select()

# Sheet1 becomes the active sheet/tab, with N9 selected, matching the target view state
$ws1.Activate()
$ws1.Range("N9").Select()
